$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 97 <- data from row 98 (B id 7323253)
$ws.Range("B97").Value = 7323253
$ws.Range("C97").Value = 'Chile Primera Division'
$ws.Range("D97").Value = 45242.83333333334
$ws.Range("E97").Value = 'Union Espanola'
$ws.Range("F97").Value = 'OHiggins'
$ws.Range("G97").Value = 3
$ws.Range("H97").Value = 3
$ws.Range("I97").Value = 1
$ws.Range("J97").Value = 2
$ws.Range("K97").Value = 'D'
$ws.Range("L97").Value = 2
$ws.Range("M97").Value = 3.4
$ws.Range("N97").Value = 3.5
$ws.Range("O97").Value = 2.1
$ws.Range("P97").Value = 3.5
$ws.Range("Q97").Value = 3.75
$ws.Range("R97").Value = -0.5
$ws.Range("S97").Value = 2.025
$ws.Range("T97").Value = 1.775
$ws.Range("U97").Value = 2.5
$ws.Range("V97").Value = 1.95
$ws.Range("W97").Value = 1.85
$ws.Range("X97").Value = -1
$ws.Range("Y97").Value = 2.5
$ws.Range("Z97").Value = -1
$ws.Range("AA97").Value = -1
$ws.Range("AB97").Value = 0.7749999999999999
$ws.Range("AC97").Value = 0.95
$ws.Range("AD97").Value = -1

# Row 98 <- data from row 97 (B id 7323186)
$ws.Range("B98").Value = 7323186
$ws.Range("C98").Value = 'Chile Primera Division'
$ws.Range("D98").Value = 45242.83333333334
$ws.Range("E98").Value = 'Coquimbo Unido'
$ws.Range("F98").Value = 'Deportes Copiapo'
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 'H'
$ws.Range("L98").Value = 2
$ws.Range("M98").Value = 3.4
$ws.Range("N98").Value = 3.5
$ws.Range("O98").Value = 1.727
$ws.Range("P98").Value = 3.8
$ws.Range("Q98").Value = 4.75
$ws.Range("R98").Value = -0.75
$ws.Range("S98").Value = 1.9
$ws.Range("T98").Value = 1.9
$ws.Range("U98").Value = 2.75
$ws.Range("V98").Value = 1.85
$ws.Range("W98").Value = 1.95
$ws.Range("X98").Value = 0.7270000000000001
$ws.Range("Y98").Value = -1
$ws.Range("Z98").Value = -1
$ws.Range("AA98").Value = 0.45
$ws.Range("AB98").Value = -0.5
$ws.Range("AC98").Value = -1
$ws.Range("AD98").Value = 0.95

# Row 102 <- data from row 103 (B id 7494646)
$ws.Range("B102").Value = 7494646
$ws.Range("C102").Value = 'Chile Primera Division'
$ws.Range("D102").Value = 45255.75
$ws.Range("E102").Value = 'OHiggins'
$ws.Range("F102").Value = 'Cobresal'
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 'D'
$ws.Range("L102").Value = 3
$ws.Range("M102").Value = 3.4
$ws.Range("N102").Value = 2.3
$ws.Range("O102").Value = 2.1
$ws.Range("P102").Value = 3.5
$ws.Range("Q102").Value = 3.5
$ws.Range("R102").Value = -0.25
$ws.Range("S102").Value = 1.8
$ws.Range("T102").Value = 2.05
$ws.Range("U102").Value = 2.75
$ws.Range("V102").Value = 1.975
$ws.Range("W102").Value = 1.875
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = 2.5
$ws.Range("Z102").Value = -1
$ws.Range("AA102").Value = -0.5
$ws.Range("AB102").Value = 0.5249999999999999
$ws.Range("AC102").Value = -1
$ws.Range("AD102").Value = 0.875

# Row 103 <- data from row 102 (B id 7494647)
$ws.Range("B103").Value = 7494647
$ws.Range("C103").Value = 'Chile Primera Division'
$ws.Range("D103").Value = 45255.75
$ws.Range("E103").Value = 'Huachipato'
$ws.Range("F103").Value = 'Universidad Catolica'
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 1
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1
$ws.Range("K103").Value = 'D'
$ws.Range("L103").Value = 2.2
$ws.Range("M103").Value = 3.4
$ws.Range("N103").Value = 3.2
$ws.Range("O103").Value = 1.8
$ws.Range("P103").Value = 3.6
$ws.Range("Q103").Value = 4.333
$ws.Range("R103").Value = -0.75
$ws.Range("S103").Value = 1.975
$ws.Range("T103").Value = 1.875
$ws.Range("U103").Value = 2.75
$ws.Range("V103").Value = 1.975
$ws.Range("W103").Value = 1.875
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = 2.6
$ws.Range("Z103").Value = -1
$ws.Range("AA103").Value = -1
$ws.Range("AB103").Value = 0.875
$ws.Range("AC103").Value = -1
$ws.Range("AD103").Value = 0.875

# Row 105 <- data from row 106 (B id 6077497)
$ws.Range("B105").Value = 6077497
$ws.Range("C105").Value = 'Chile Primera Division'
$ws.Range("D105").Value = 45256.85416666666
$ws.Range("E105").Value = 'Deportes Copiapo'
$ws.Range("F105").Value = 'Nublense'
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 1
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 'D'
$ws.Range("L105").Value = 2.6
$ws.Range("M105").Value = 3.4
$ws.Range("N105").Value = 2.6
$ws.Range("O105").Value = 2.8
$ws.Range("P105").Value = 3.2
$ws.Range("Q105").Value = 2.7
$ws.Range("R105").Value = 0
$ws.Range("S105").Value = 1.95
$ws.Range("T105").Value = 1.9
$ws.Range("U105").Value = 2.25
$ws.Range("V105").Value = 2
$ws.Range("W105").Value = 1.85
$ws.Range("X105").Value = -1
$ws.Range("Y105").Value = 2.2
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = 0
$ws.Range("AB105").Value = 0
$ws.Range("AC105").Value = -0.5
$ws.Range("AD105").Value = 0.425

# Row 106 <- data from row 105 (B id 6077763)
$ws.Range("B106").Value = 6077763
$ws.Range("C106").Value = 'Chile Primera Division'
$ws.Range("D106").Value = 45256.85416666666
$ws.Range("E106").Value = 'Curico Unido'
$ws.Range("F106").Value = 'Magallanes'
$ws.Range("G106").Value = 3
$ws.Range("H106").Value = 4
$ws.Range("I106").Value = 2
$ws.Range("J106").Value = 2
$ws.Range("K106").Value = 'A'
$ws.Range("L106").Value = 2.15
$ws.Range("M106").Value = 3.5
$ws.Range("N106").Value = 3.2
$ws.Range("O106").Value = 2.625
$ws.Range("P106").Value = 3.5
$ws.Range("Q106").Value = 2.6
$ws.Range("R106").Value = 0
$ws.Range("S106").Value = 1.95
$ws.Range("T106").Value = 1.9
$ws.Range("U106").Value = 2.75
$ws.Range("V106").Value = 1.975
$ws.Range("W106").Value = 1.875
$ws.Range("X106").Value = -1
$ws.Range("Y106").Value = -1
$ws.Range("Z106").Value = 1.6
$ws.Range("AA106").Value = -1
$ws.Range("AB106").Value = 0.8999999999999999
$ws.Range("AC106").Value = 0.9750000000000001
$ws.Range("AD106").Value = -1

# Row 109 <- data from row 110 (B id 6077498)
$ws.Range("B109").Value = 6077498
$ws.Range("C109").Value = 'Chile Primera Division'
$ws.Range("D109").Value = 45262.75
$ws.Range("E109").Value = 'Universidad Catolica'
$ws.Range("F109").Value = 'Deportes Copiapo'
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 2
$ws.Range("I109").Value = 2
$ws.Range("J109").Value = 2
$ws.Range("K109").Value = 'D'
$ws.Range("L109").Value = 1.65
$ws.Range("M109").Value = 3.8
$ws.Range("N109").Value = 5.25
$ws.Range("O109").Value = 1.909
$ws.Range("P109").Value = 3.6
$ws.Range("Q109").Value = 4.2
$ws.Range("R109").Value = -0.5
$ws.Range("S109").Value = 1.85
$ws.Range("T109").Value = 2
$ws.Range("U109").Value = 2.75
$ws.Range("V109").Value = 2.025
$ws.Range("W109").Value = 1.825
$ws.Range("X109").Value = -1
$ws.Range("Y109").Value = 2.6
$ws.Range("Z109").Value = -1
$ws.Range("AA109").Value = -1
$ws.Range("AB109").Value = 1
$ws.Range("AC109").Value = 1.025
$ws.Range("AD109").Value = -1

# Row 110 <- data from row 109 (B id 6078266)
$ws.Range("B110").Value = 6078266
$ws.Range("C110").Value = 'Chile Primera Division'
$ws.Range("D110").Value = 45262.75
$ws.Range("E110").Value = 'Palestino'
$ws.Range("F110").Value = 'Curico Unido'
$ws.Range("G110").Value = 4
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 'H'
$ws.Range("L110").Value = 1.533
$ws.Range("M110").Value = 4
$ws.Range("N110").Value = 6
$ws.Range("O110").Value = 1.363
$ws.Range("P110").Value = 4.75
$ws.Range("Q110").Value = 7.5
$ws.Range("R110").Value = -1.5
$ws.Range("S110").Value = 2.025
$ws.Range("T110").Value = 1.825
$ws.Range("U110").Value = 3
$ws.Range("V110").Value = 1.9
$ws.Range("W110").Value = 1.95
$ws.Range("X110").Value = 0.363
$ws.Range("Y110").Value = -1
$ws.Range("Z110").Value = -1
$ws.Range("AA110").Value = 1.025
$ws.Range("AB110").Value = -1
$ws.Range("AC110").Value = 0.8999999999999999
$ws.Range("AD110").Value = -1

# Row 112 <- data from row 114 (B id 6077767)
$ws.Range("B112").Value = 6077767
$ws.Range("C112").Value = 'Chile Primera Division'
$ws.Range("D112").Value = 45263.75
$ws.Range("E112").Value = 'Nublense'
$ws.Range("F112").Value = 'Huachipato'
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 1
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 'A'
$ws.Range("L112").Value = 2.75
$ws.Range("M112").Value = 3.4
$ws.Range("N112").Value = 2.45
$ws.Range("O112").Value = 2.875
$ws.Range("P112").Value = 3.3
$ws.Range("Q112").Value = 2.5
$ws.Range("R112").Value = 0
$ws.Range("S112").Value = 2.05
$ws.Range("T112").Value = 1.8
$ws.Range("U112").Value = 2.25
$ws.Range("V112").Value = 1.8
$ws.Range("W112").Value = 2.05
$ws.Range("X112").Value = -1
$ws.Range("Y112").Value = -1
$ws.Range("Z112").Value = 1.5
$ws.Range("AA112").Value = -1
$ws.Range("AB112").Value = 0.8
$ws.Range("AC112").Value = -1
$ws.Range("AD112").Value = 1.05

# Row 113 <- data from row 112 (B id 6078996)
$ws.Range("B113").Value = 6078996
$ws.Range("C113").Value = 'Chile Primera Division'
$ws.Range("D113").Value = 45263.75
$ws.Range("E113").Value = 'Colo Colo'
$ws.Range("F113").Value = 'Union Espanola'
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 2
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 'A'
$ws.Range("L113").Value = 1.4
$ws.Range("M113").Value = 4.333
$ws.Range("N113").Value = 7
$ws.Range("O113").Value = 1.285
$ws.Range("P113").Value = 5.5
$ws.Range("Q113").Value = 11
$ws.Range("R113").Value = -1.5
$ws.Range("S113").Value = 1.9
$ws.Range("T113").Value = 1.95
$ws.Range("U113").Value = 3
$ws.Range("V113").Value = 2
$ws.Range("W113").Value = 1.85
$ws.Range("X113").Value = -1
$ws.Range("Y113").Value = -1
$ws.Range("Z113").Value = 10
$ws.Range("AA113").Value = -1
$ws.Range("AB113").Value = 0.95
$ws.Range("AC113").Value = -1
$ws.Range("AD113").Value = 0.8500000000000001

# Row 114 <- data from row 113 (B id 6078263)
$ws.Range("B114").Value = 6078263
$ws.Range("C114").Value = 'Chile Primera Division'
$ws.Range("D114").Value = 45263.75
$ws.Range("E114").Value = 'Cobresal'
$ws.Range("F114").Value = 'Universidad de Chile'
$ws.Range("G114").Value = 4
$ws.Range("H114").Value = 3
$ws.Range("I114").Value = 1
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 'H'
$ws.Range("L114").Value = 2.05
$ws.Range("M114").Value = 3.4
$ws.Range("N114").Value = 3.5
$ws.Range("O114").Value = 2.05
$ws.Range("P114").Value = 3.6
$ws.Range("Q114").Value = 3.5
$ws.Range("R114").Value = -0.5
$ws.Range("S114").Value = 2
$ws.Range("T114").Value = 1.8
$ws.Range("U114").Value = 2.75
$ws.Range("V114").Value = 1.9
$ws.Range("W114").Value = 1.9
$ws.Range("X114").Value = 1.05
$ws.Range("Y114").Value = -1
$ws.Range("Z114").Value = -1
$ws.Range("AA114").Value = 1
$ws.Range("AB114").Value = -1
$ws.Range("AC114").Value = 0.8999999999999999
$ws.Range("AD114").Value = -1

# Row 115 <- data from row 117 (B id 6143704)
$ws.Range("B115").Value = 6143704
$ws.Range("C115").Value = 'Chile Primera Division'
$ws.Range("D115").Value = 45268.75
$ws.Range("E115").Value = 'Curico Unido'
$ws.Range("F115").Value = 'Colo Colo'
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 1
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 'A'
$ws.Range("L115").Value = 6.5
$ws.Range("M115").Value = 4.75
$ws.Range("N115").Value = 1.4
$ws.Range("O115").Value = 12
$ws.Range("P115").Value = 8.5
$ws.Range("Q115").Value = 1.166
$ws.Range("R115").Value = 2
$ws.Range("S115").Value = 2
$ws.Range("T115").Value = 1.8
$ws.Range("U115").Value = 3.25
$ws.Range("V115").Value = 1.875
$ws.Range("W115").Value = 1.925
$ws.Range("X115").Value = -1
$ws.Range("Y115").Value = -1
$ws.Range("Z115").Value = 0.1659999999999999
$ws.Range("AA115").Value = 1
$ws.Range("AB115").Value = -1
$ws.Range("AC115").Value = -1
$ws.Range("AD115").Value = 0.925

# Row 117 <- data from row 115 (B id 6078997)
$ws.Range("B117").Value = 6078997
$ws.Range("C117").Value = 'Chile Primera Division'
$ws.Range("D117").Value = 45268.75
$ws.Range("E117").Value = 'Union Espanola'
$ws.Range("F117").Value = 'Cobresal'
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 'H'
$ws.Range("L117").Value = 3.8
$ws.Range("M117").Value = 3.6
$ws.Range("N117").Value = 1.909
$ws.Range("O117").Value = 2.7
$ws.Range("P117").Value = 3.6
$ws.Range("Q117").Value = 2.45
$ws.Range("R117").Value = 0
$ws.Range("S117").Value = 1.975
$ws.Range("T117").Value = 1.825
$ws.Range("U117").Value = 2.75
$ws.Range("V117").Value = 1.775
$ws.Range("W117").Value = 2.025
$ws.Range("X117").Value = 1.7
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = -1
$ws.Range("AA117").Value = 0.9750000000000001
$ws.Range("AB117").Value = -1
$ws.Range("AC117").Value = -1
$ws.Range("AD117").Value = 1.025

# Row 118 <- data from row 119 (B id 6077768)
$ws.Range("B118").Value = 6077768
$ws.Range("C118").Value = 'Chile Primera Division'
$ws.Range("D118").Value = 45269.75
$ws.Range("E118").Value = 'Union La Calera'
$ws.Range("F118").Value = 'Universidad Catolica'
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 3
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 1
$ws.Range("K118").Value = 'A'
$ws.Range("L118").Value = 2.05
$ws.Range("M118").Value = 3.5
$ws.Range("N118").Value = 3.4
$ws.Range("O118").Value = 2.05
$ws.Range("P118").Value = 3.6
$ws.Range("Q118").Value = 3.4
$ws.Range("R118").Value = -0.25
$ws.Range("S118").Value = 1.8
$ws.Range("T118").Value = 2
$ws.Range("U118").Value = 2.75
$ws.Range("V118").Value = 1.975
$ws.Range("W118").Value = 1.825
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 2.4
$ws.Range("AA118").Value = -1
$ws.Range("AB118").Value = 1
$ws.Range("AC118").Value = 0.4875
$ws.Range("AD118").Value = -0.5

# Row 119 <- data from row 120 (B id 6078268)
$ws.Range("B119").Value = 6078268
$ws.Range("C119").Value = 'Chile Primera Division'
$ws.Range("D119").Value = 45269.75
$ws.Range("E119").Value = 'OHiggins'
$ws.Range("F119").Value = 'Palestino'
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 1
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 1
$ws.Range("K119").Value = 'A'
$ws.Range("L119").Value = 3.1
$ws.Range("M119").Value = 3.3
$ws.Range("N119").Value = 2.3
$ws.Range("O119").Value = 2.9
$ws.Range("P119").Value = 3.4
$ws.Range("Q119").Value = 2.375
$ws.Range("R119").Value = 0.25
$ws.Range("S119").Value = 1.8
$ws.Range("T119").Value = 2
$ws.Range("U119").Value = 2.75
$ws.Range("V119").Value = 2
$ws.Range("W119").Value = 1.8
$ws.Range("X119").Value = -1
$ws.Range("Y119").Value = -1
$ws.Range("Z119").Value = 1.375
$ws.Range("AA119").Value = -1
$ws.Range("AB119").Value = 1
$ws.Range("AC119").Value = -1
$ws.Range("AD119").Value = 0.8

# Row 120 <- data from row 118 (B id 6077499)
$ws.Range("B120").Value = 6077499
$ws.Range("C120").Value = 'Chile Primera Division'
$ws.Range("D120").Value = 45269.75
$ws.Range("E120").Value = 'Deportes Copiapo'
$ws.Range("F120").Value = 'Everton de Vina'
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 2
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 'H'
$ws.Range("L120").Value = 2.1
$ws.Range("M120").Value = 3.4
$ws.Range("N120").Value = 3.4
$ws.Range("O120").Value = 2.9
$ws.Range("P120").Value = 3.4
$ws.Range("Q120").Value = 2.4
$ws.Range("R120").Value = 0.25
$ws.Range("S120").Value = 1.775
$ws.Range("T120").Value = 2.1
$ws.Range("U120").Value = 2.75
$ws.Range("V120").Value = 1.85
$ws.Range("W120").Value = 2
$ws.Range("X120").Value = 1.9
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = -1
$ws.Range("AA120").Value = 0.7749999999999999
$ws.Range("AB120").Value = -1
$ws.Range("AC120").Value = -1
$ws.Range("AD120").Value = 1

# Row 121 <- data from row 122 (B id 6078998)
$ws.Range("B121").Value = 6078998
$ws.Range("C121").Value = 'Chile Primera Division'
$ws.Range("D121").Value = 45269.75
$ws.Range("E121").Value = 'Magallanes'
$ws.Range("F121").Value = 'Coquimbo Unido'
$ws.Range("G121").Value = 2
$ws.Range("H121").Value = 3
$ws.Range("I121").Value = 1
$ws.Range("J121").Value = 1
$ws.Range("K121").Value = 'A'
$ws.Range("L121").Value = 1.909
$ws.Range("M121").Value = 3.6
$ws.Range("N121").Value = 3.8
$ws.Range("O121").Value = 2.15
$ws.Range("P121").Value = 3.75
$ws.Range("Q121").Value = 3.1
$ws.Range("R121").Value = -0.25
$ws.Range("S121").Value = 1.85
$ws.Range("T121").Value = 1.95
$ws.Range("U121").Value = 3
$ws.Range("V121").Value = 1.85
$ws.Range("W121").Value = 1.95
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = -1
$ws.Range("Z121").Value = 2.1
$ws.Range("AA121").Value = -1
$ws.Range("AB121").Value = 0.95
$ws.Range("AC121").Value = 0.8500000000000001
$ws.Range("AD121").Value = -1

# Row 122 <- data from row 121 (B id 6078269)
$ws.Range("B122").Value = 6078269
$ws.Range("C122").Value = 'Chile Primera Division'
$ws.Range("D122").Value = 45269.75
$ws.Range("E122").Value = 'Universidad de Chile'
$ws.Range("F122").Value = 'Nublense'
$ws.Range("G122").Value = 3
$ws.Range("H122").Value = 1
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = 1
$ws.Range("K122").Value = 'H'
$ws.Range("L122").Value = 1.85
$ws.Range("M122").Value = 3.4
$ws.Range("N122").Value = 4.333
$ws.Range("O122").Value = 1.8
$ws.Range("P122").Value = 3.6
$ws.Range("Q122").Value = 4.5
$ws.Range("R122").Value = -0.75
$ws.Range("S122").Value = 1.925
$ws.Range("T122").Value = 1.925
$ws.Range("U122").Value = 2.5
$ws.Range("V122").Value = 2.025
$ws.Range("W122").Value = 1.825
$ws.Range("X122").Value = 0.8
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = -1
$ws.Range("AA122").Value = 0.925
$ws.Range("AB122").Value = -1
$ws.Range("AC122").Value = 1.025
$ws.Range("AD122").Value = -1

# Row 139 <- data from row 140 (B id 7723532)
$ws.Range("B139").Value = 7723532
$ws.Range("C139").Value = 'Chile Primera Division'
$ws.Range("D139").Value = 45353.85416666666
$ws.Range("E139").Value = 'Nublense'
$ws.Range("F139").Value = 'Audax Italiano'
$ws.Range("G139").Value = 1
$ws.Range("H139").Value = 2
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 1
$ws.Range("K139").Value = 'A'
$ws.Range("L139").Value = 2
$ws.Range("M139").Value = 3.25
$ws.Range("N139").Value = 3.6
$ws.Range("O139").Value = 2.25
$ws.Range("P139").Value = 3.2
$ws.Range("Q139").Value = 3.5
$ws.Range("R139").Value = -0.25
$ws.Range("S139").Value = 1.875
$ws.Range("T139").Value = 1.925
$ws.Range("U139").Value = 2.25
$ws.Range("V139").Value = 1.9
$ws.Range("W139").Value = 1.9
$ws.Range("X139").Value = -1
$ws.Range("Y139").Value = -1
$ws.Range("Z139").Value = 2.5
$ws.Range("AA139").Value = -1
$ws.Range("AB139").Value = 0.925
$ws.Range("AC139").Value = 0.8999999999999999
$ws.Range("AD139").Value = -1

# Row 140 <- data from row 139 (B id 7723531)
$ws.Range("B140").Value = 7723531
$ws.Range("C140").Value = 'Chile Primera Division'
$ws.Range("D140").Value = 45353.85416666666
$ws.Range("E140").Value = 'Union Espanola'
$ws.Range("F140").Value = 'Coquimbo Unido'
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 1
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 'H'
$ws.Range("L140").Value = 2.05
$ws.Range("M140").Value = 3.4
$ws.Range("N140").Value = 3.25
$ws.Range("O140").Value = 2.05
$ws.Range("P140").Value = 3.5
$ws.Range("Q140").Value = 3.5
$ws.Range("R140").Value = -0.25
$ws.Range("S140").Value = 1.775
$ws.Range("T140").Value = 2.025
$ws.Range("U140").Value = 2.5
$ws.Range("V140").Value = 1.8
$ws.Range("W140").Value = 2
$ws.Range("X140").Value = 1.05
$ws.Range("Y140").Value = -1
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 0.7749999999999999
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = -1
$ws.Range("AD140").Value = 1

# Row 179 <- data from row 180 (B id 7723557)
$ws.Range("B179").Value = 7723557
$ws.Range("C179").Value = 'Chile Primera Division'
$ws.Range("D179").Value = 45396.5625
$ws.Range("E179").Value = 'Universidad de Chile'
$ws.Range("F179").Value = 'Coquimbo Unido'
$ws.Range("G179").Value = 1
$ws.Range("H179").Value = 1
$ws.Range("I179").Value = 0
$ws.Range("J179").Value = 0
$ws.Range("K179").Value = 'D'
$ws.Range("L179").Value = 1.8
$ws.Range("M179").Value = 3.6
$ws.Range("N179").Value = 4.333
$ws.Range("O179").Value = 1.95
$ws.Range("P179").Value = 3.5
$ws.Range("Q179").Value = 4
$ws.Range("R179").Value = -0.5
$ws.Range("S179").Value = 1.925
$ws.Range("T179").Value = 1.925
$ws.Range("U179").Value = 2.5
$ws.Range("V179").Value = 1.975
$ws.Range("W179").Value = 1.875
$ws.Range("X179").Value = -1
$ws.Range("Y179").Value = 2.5
$ws.Range("Z179").Value = -1
$ws.Range("AA179").Value = -1
$ws.Range("AB179").Value = 0.925
$ws.Range("AC179").Value = -1
$ws.Range("AD179").Value = 0.875

# Row 180 <- data from row 179 (B id 7723561)
$ws.Range("B180").Value = 7723561
$ws.Range("C180").Value = 'Chile Primera Division'
$ws.Range("D180").Value = 45396.5625
$ws.Range("E180").Value = 'OHiggins'
$ws.Range("F180").Value = 'Huachipato'
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 0
$ws.Range("I180").Value = 0
$ws.Range("J180").Value = 0
$ws.Range("K180").Value = 'D'
$ws.Range("L180").Value = 2.4
$ws.Range("M180").Value = 3.2
$ws.Range("N180").Value = 3
$ws.Range("O180").Value = 2.45
$ws.Range("P180").Value = 3.2
$ws.Range("Q180").Value = 3
$ws.Range("R180").Value = -0.25
$ws.Range("S180").Value = 2.025
$ws.Range("T180").Value = 1.775
$ws.Range("U180").Value = 2.25
$ws.Range("V180").Value = 1.825
$ws.Range("W180").Value = 1.975
$ws.Range("X180").Value = -1
$ws.Range("Y180").Value = 2.2
$ws.Range("Z180").Value = -1
$ws.Range("AA180").Value = -0.5
$ws.Range("AB180").Value = 0.3875
$ws.Range("AC180").Value = -1
$ws.Range("AD180").Value = 0.9750000000000001

# Row 231 standalone value updates
$ws.Range("V231").Value = 1.95
$ws.Range("W231").Value = 1.9
